$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18 (QSY_B_072): fill in the scar_number that was missing ---
$ws.Range("B18").Value = 4

# The SDI column (G = scar_number / surface area) uses a shared formula
# across G9:G21; row 18 now has the data it needs, so extend it down.
$ws.Range("G18").Formula = "=B18/D18"
# Writing the formula pulls in neighbouring-cell formatting; put G18 back
# to an unformatted cell like the rest of the G column (G9:G17, G21).
$ws.Range("G18").ClearFormats()

# A18 ("QSY_B_072") was styled differently (plain black) while its data
# was incomplete; now that scar_number is filled in, match the red font
# used by the other completed rows (A10:A17, A5, A6, A9, A21, ...).
$ws.Range("A18").Font.Color = 255

# Reflect where the user's cursor ended up after these edits.
$ws.Range("C22").Select()
